# Update TestActions / TestBase test_suite sheet to reflect the new site.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_suite")

# Rows 15 (View_Map) and 18 (ResetPassword_viaForgotPassword) now expect "N"
# instead of "Y" in column B.
$ws.Range("B15").Value = "N"
$ws.Range("B18").Value = "N"

# Move the active selection to B15 (was B19).
$ws.Range("B15").Select()
